$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 12.1382
$ws.Range("B8").Value = 4.824999999999999
$ws.Range("B10").Value = 8.660000000000004
$ws.Range("B12").Value = 6.867500000000002
$ws.Range("D13").Value = -7.4968
$ws.Range("B18").Value = 6.431400000000001
$ws.Range("E20").Value = 13.1987
